# Generate Report for Handback
# Updates the Correspond Handoff/Handback timestamps for the
# "79ffca11-..." row on the zh-cn and de-de sheets, and refreshes the
# "Latest HO Xliff Generate Date" summary on the Overview sheet to match.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-25 16:50:55"
$zhcn.Range("K2").Value = "2016-08-25 16:51:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-25 16:50:59"
$dede.Range("K2").Value = "2016-08-25 16:51:25"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-25 16:50:59"
